$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 61) to the FWHM data table, for the FSR data run of
# sg_rr_100_028 2023-12-08 16-58-05
$ws.Range("A61").Value = "sg_rr_100_028 2023-12-08 16-58-05.csv"
$ws.Range("B61").Value = 0.01
$ws.Range("C61").Value = 1000
$ws.Range("D61").Value = 5001
$ws.Range("E61").Value = 1530
$ws.Range("F61").Value = 1570
$ws.Range("G61").Value = 0.5
$ws.Range("H61").Value = "(approx_fsr/2)/wavelength step size"
$ws.Range("I61").Value = 1
$ws.Range("J61").Value = 0.98274999999999801
$ws.Range("K61").Value = 0.0032814064370514399
$ws.Range("L61").Value = "yes"
$ws.Range("M61").Value = 0.102827046790518
$ws.Range("N61").Value = 0.00213794013833199
$ws.Range("O61").Value = "reduced approx fsr a bit, to see if this had any affect on fsr calculation as above, half the approx fsr was quite close to actual calculated fsr."

# Update the view to reflect where the user ended up after the edit
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Range("O60:O61").Select()
